$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab (workbook.xml <sheet name="...">)
$ws.Name = "BetaFiberA"

# Add new row 16, mirroring row 15's formatting (bold/border style on col A)
# by copying the existing row's formats down, then filling in the new values.
$ws.Range("A15:P15").Copy()
$ws.Range("A16:P16").PasteSpecial(-4122)

$ws.Range("A16").Value = 14
$ws.Range("B16").Value2 = $ws.Range("B15").Value2

for ($col = 3; $col -le 16; $col++) {
    $ws.Cells.Item(16, $col).Value = 1
}
